# Correcciones en script de carga y archivos de prueba para actividades de ETL
# Ventas sheet: add an id_pedido column (D) and add a new sale row (Audifonos -> row 7),
# re-deriving the product/price values for rows 4-6 to reflect order groupings.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ventas")

# New header for column D (copy the header formatting from an existing header cell)
$ws.Cells.Item(1, 4).Value = "id_pedido"
$ws.Range("C1").Copy()
$ws.Range("D1").PasteSpecial(-4122)  # xlPasteFormats

# Rebuild rows 2-7 with the updated producto/precio/id_pedido data
$ws.Cells.Item(2, 1).Value = 1
$ws.Cells.Item(2, 2).Value = "Laptop"
$ws.Cells.Item(2, 3).Value = 1200
$ws.Cells.Item(2, 4).Value = 1

$ws.Cells.Item(3, 1).Value = 2
$ws.Cells.Item(3, 2).Value = "Mouse"
$ws.Cells.Item(3, 3).Value = 25
$ws.Cells.Item(3, 4).Value = 1

$ws.Cells.Item(4, 1).Value = 3
$ws.Cells.Item(4, 2).Value = "Mouse"
$ws.Cells.Item(4, 3).Value = 25
$ws.Cells.Item(4, 4).Value = 2

$ws.Cells.Item(5, 1).Value = 4
$ws.Cells.Item(5, 2).Value = "Teclado"
$ws.Cells.Item(5, 3).Value = 80
$ws.Cells.Item(5, 4).Value = 3

$ws.Cells.Item(6, 1).Value = 5
$ws.Cells.Item(6, 2).Value = "Monitor"
$ws.Cells.Item(6, 3).Value = 300
$ws.Cells.Item(6, 4).Value = 3

$ws.Cells.Item(7, 1).Value = 6
$ws.Cells.Item(7, 2).Value = "Audífonos"
$ws.Cells.Item(7, 3).Value = 150
$ws.Cells.Item(7, 4).Value = 4
